# Remove the stray "You might also like" text that was accidentally scraped
# into the lyrics column (column C) for every data row of the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)   # column C
    $val = $cell.Value2

    if ($null -ne $val -and $val -is [string] -and $val.Contains("You might also like")) {
        $new = $val -replace ' You might also like ', ' '
        $new = $new -replace ' You might also like$', ''
        $new = $new -replace '^You might also like ', ''
        $new = $new -replace 'You might also like', ''

        $cell.Value2 = $new
    }
}
